$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 322, pushing the existing row 322 (and all
# rows below it) down by one. Excel copies the formatting of the row above
# into the freshly inserted row (matches the date-formatted style on column D).
$ws.Rows("322:322").Insert()

# Populate the newly inserted row 322 with the new daily price record.
$ws.Cells.Item(322, 1).Value = 8
$ws.Cells.Item(322, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(322, 3).Value = "Coquimbo"
$ws.Cells.Item(322, 4).Value = 44876
$ws.Cells.Item(322, 5).Value = 4
$ws.Cells.Item(322, 6).Value = 100114013
$ws.Cells.Item(322, 7).Value = "Zanahoria"
$ws.Cells.Item(322, 8).Value = "Sin especificar"
$ws.Cells.Item(322, 9).Value = "Primera"
$ws.Cells.Item(322, 10).Value = 400
$ws.Cells.Item(322, 11).Value = 14000
$ws.Cells.Item(322, 12).Value = 15000
$ws.Cells.Item(322, 13).Value = 14500
$ws.Cells.Item(322, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(322, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(322, 16).Value = 725
$ws.Cells.Item(322, 17).Value = 20
$ws.Cells.Item(322, 18).Value = "Hortaliza"
